# Add a new weekly price record for "Zapallo italiano" (Vega Monumental
# Concepción) as row 89, pushing the existing rows 89-110 down to 90-111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 89 (shifts rows 89..110 -> 90..111)
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row with the new observation
$ws.Cells.Item(89, 1).Value = 11
$ws.Cells.Item(89, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(89, 3).Value = "Bíobío"
$ws.Cells.Item(89, 4).Value = 44642
$ws.Cells.Item(89, 5).Value = 8
$ws.Cells.Item(89, 6).Value = 100112032
$ws.Cells.Item(89, 7).Value = "Zapallo italiano"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 220
$ws.Cells.Item(89, 11).Value = 13000
$ws.Cells.Item(89, 12).Value = 14000
$ws.Cells.Item(89, 13).Value = 13455
$ws.Cells.Item(89, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(89, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value = 224
$ws.Cells.Item(89, 17).Value = 60
$ws.Cells.Item(89, 18).Value = "Hortaliza"
